$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 08:54:55"
$wsZhCn.Range("H2").Value = "2016-03-19 08:55:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 08:54:58"
$wsDeDe.Range("H2").Value = "2016-03-19 08:55:19"
